$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price column D, Volume(1h) column E).
# Values are written as text (NumberFormat "@") to preserve exact
# formatting (trailing zeros, percent signs) exactly like the source data.
$updates = @(
    @{ Cell = "D2"; Value = "261.36" }
    @{ Cell = "D3"; Value = "27.08" }
    @{ Cell = "E3"; Value = "0.65%" }
    @{ Cell = "D4"; Value = "4.714" }
    @{ Cell = "E4"; Value = "0.48%" }
    @{ Cell = "D5"; Value = "0.06207" }
    @{ Cell = "E5"; Value = "2.74%" }
    @{ Cell = "D6"; Value = "6.726" }
    @{ Cell = "E6"; Value = "0.78%" }
    @{ Cell = "D7"; Value = "0.8503" }
    @{ Cell = "E7"; Value = "-1.06%" }
    @{ Cell = "D8"; Value = "0.9112" }
    @{ Cell = "E8"; Value = "-0.78%" }
    @{ Cell = "D9"; Value = "0.1406" }
    @{ Cell = "E9"; Value = "0.75%" }
    @{ Cell = "D10"; Value = "0.04704" }
    @{ Cell = "E10"; Value = "-12.47%" }
    @{ Cell = "D11"; Value = "0.07094" }
    @{ Cell = "E11"; Value = "0.18%" }
    @{ Cell = "D12"; Value = "0.03168" }
    @{ Cell = "E12"; Value = "3.20%" }
    @{ Cell = "E13"; Value = "-0.75%" }
    @{ Cell = "D14"; Value = "0.001535" }
    @{ Cell = "E14"; Value = "0.37%" }
    @{ Cell = "D15"; Value = "0.0006178" }
    @{ Cell = "E15"; Value = "1.68%" }
    @{ Cell = "D16"; Value = "0.006058" }
    @{ Cell = "E16"; Value = "-0.17%" }
    @{ Cell = "D17"; Value = "3.464" }
    @{ Cell = "E17"; Value = "-0.10%" }
    @{ Cell = "E18"; Value = "-0.14%" }
    @{ Cell = "E19"; Value = "0.57%" }
    @{ Cell = "E20"; Value = "-0.68%" }
    @{ Cell = "D21"; Value = "0.1310" }
    @{ Cell = "E21"; Value = "0.95%" }
    @{ Cell = "D22"; Value = "4.099" }
    @{ Cell = "E22"; Value = "-0.85%" }
    @{ Cell = "D23"; Value = "0.04226" }
    @{ Cell = "E23"; Value = "-0.45%" }
    @{ Cell = "D24"; Value = "0.001218" }
    @{ Cell = "E24"; Value = "0.19%" }
    @{ Cell = "D25"; Value = "0.004128" }
    @{ Cell = "E25"; Value = "2.57%" }
    @{ Cell = "E26"; Value = "0.17%" }
    @{ Cell = "D40"; Value = "0.03915" }
    @{ Cell = "E40"; Value = "1.71%" }
    @{ Cell = "D41"; Value = "0.1112" }
    @{ Cell = "E41"; Value = "-0.22%" }
    @{ Cell = "D42"; Value = "0.004130" }
    @{ Cell = "E42"; Value = "2.62%" }
    @{ Cell = "E43"; Value = "-0.65%" }
    @{ Cell = "D44"; Value = "0.01353" }
    @{ Cell = "E44"; Value = "-10.60%" }
    @{ Cell = "E45"; Value = "0.32%" }
    @{ Cell = "E46"; Value = "0.17%" }
    @{ Cell = "E47"; Value = "-34.11%" }
    @{ Cell = "D48"; Value = "0.1700" }
    @{ Cell = "E48"; Value = "28.70%" }
    @{ Cell = "E49"; Value = "0.17%" }
    @{ Cell = "E50"; Value = "0.17%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
